# Insert a new weekly price record for "Apio" (Macroferia Regional de Talca)
# right before the existing row 166, shifting rows 166:295 down to 167:296
# (dimension grows from A1:R295 to A1:R296).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 166 - this pushes the old row166..row295 down to
# row167..row296, carrying their formatting (incl. the date style on column D).
$ws.Rows.Item(166).Insert()

# Populate the newly inserted row 166 with the new weekly data point.
$ws.Cells.Item(166, 1).Value  = 5
$ws.Cells.Item(166, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(166, 3).Value  = "Maule"
$ws.Cells.Item(166, 4).Value  = 45068
$ws.Cells.Item(166, 5).Value  = 7
$ws.Cells.Item(166, 6).Value  = 100112017
$ws.Cells.Item(166, 7).Value  = "Apio"
$ws.Cells.Item(166, 8).Value  = "Americana (o)"
$ws.Cells.Item(166, 9).Value  = "Primera"
$ws.Cells.Item(166, 10).Value = 700
$ws.Cells.Item(166, 11).Value = 6000
$ws.Cells.Item(166, 12).Value = 6000
$ws.Cells.Item(166, 13).Value = 6000
$ws.Cells.Item(166, 14).Value = "$/docena de matas"
$ws.Cells.Item(166, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(166, 16).Value = 1000
$ws.Cells.Item(166, 17).Value = 6
$ws.Cells.Item(166, 18).Value = "Hortaliza"
